$d = $word.ActiveDocument

$replacements = @(
    @("530÷5=106, 0", "421÷8=52, 5"),
    @("215÷3=71, 2", "847÷6=141, 1"),
    @("788÷6=131, 2", "535÷8=66, 7"),
    @("380÷7=54, 2", "776÷3=258, 2"),
    @("777÷7=111, 0", "988÷8=123, 4"),
    @("122÷3=40, 2", "977÷6=162, 5"),
    @("683÷7=97, 4", "160÷8=20, 0"),
    @("271÷7=38, 5", "638÷4=159, 2"),
    @("486÷8=60, 6", "104÷8=13, 0"),
    @("329÷9=36, 5", "279÷8=34, 7"),
    @("361÷9=40, 1", "684÷4=171, 0"),
    @("162÷7=23, 1", "559÷6=93, 1"),
    @("363÷6=60, 3", "834÷3=278, 0"),
    @("654÷3=218, 0", "445÷2=222, 1"),
    @("256÷5=51, 1", "991÷9=110, 1"),
    @("852÷8=106, 4", "609÷5=121, 4"),
    @("527÷5=105, 2", "610÷3=203, 1"),
    @("243÷4=60, 3", "252÷5=50, 2"),
    @("445÷9=49, 4", "598÷7=85, 3"),
    @("669÷5=133, 4", "180÷4=45, 0"),
    @("267÷9=29, 6", "517÷9=57, 4"),
    @("801÷9=89, 0", "957÷3=319, 0"),
    @("249÷7=35, 4", "347÷4=86, 3"),
    @("776÷8=97, 0", "483÷2=241, 1"),
    @("121÷4=30, 1", "742÷6=123, 4")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
